$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 110, pushing existing rows 110-170 down to 111-171
$ws.Rows("110:110").Insert()

# Populate the newly inserted row 110 with the new data record
$ws.Cells.Item(110, 1).Value = 1
$ws.Cells.Item(110, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(110, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(110, 4).Value = 44529
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
$ws.Cells.Item(110, 5).Value = 15
$ws.Cells.Item(110, 6).Value = "Fruta"
$ws.Cells.Item(110, 7).Value = 100102
$ws.Cells.Item(110, 8).Value = "Cítricos"
$ws.Cells.Item(110, 9).Value = 100102003
$ws.Cells.Item(110, 10).Value = "Limón"
$ws.Cells.Item(110, 11).Value = "Tahití"
$ws.Cells.Item(110, 12).Value = "Primera"
$ws.Cells.Item(110, 13).Value = 200
$ws.Cells.Item(110, 14).Value = 26000
$ws.Cells.Item(110, 15).Value = 27000
$ws.Cells.Item(110, 16).Value = 26500
$ws.Cells.Item(110, 17).Value = "`$/caja 24 kilos"
$ws.Cells.Item(110, 18).Value = "Perú"
$ws.Cells.Item(110, 19).Value = 1104
$ws.Cells.Item(110, 20).Value = 24
